$d = $word.ActiveDocument

function New-ParaXml([string]$inner) {
    return '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $inner + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Insert-ParaAfter($paragraph, [string]$innerBody) {
    # Insert a brand-new <w:p> right after $paragraph, by targeting the
    # zero-length point just before its trailing paragraph mark.
    $r = $d.Range($paragraph.Range.End - 1, $paragraph.Range.End - 1)
    $xml = New-ParaXml('<w:body>' + $innerBody + '</w:body>')
    $r.InsertXML($xml)
}

function Replace-ParaContent($paragraph, [string]$innerBody) {
    # Replace a paragraph's whole content (pPr + runs) in place.
    $r = $paragraph.Range
    $xml = New-ParaXml('<w:body>' + $innerBody + '</w:body>')
    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) After P4 ("2. 종량제봉투의 용량 결정 필요"), insert new paragraph:
#    "3. 서보 모터의 적용 방식"
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$body1 = '<w:p><w:pPr><w:ind w:firstLine="195"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">3. </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>서보</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 모터의 적용 방식</w:t></w:r>' + `
    '</w:p>'
Insert-ParaAfter $p4 $body1

Write-Host "step1 done, count=$($d.Paragraphs.Count)"

# ---------------------------------------------------------------------
# 2) After "2. 보관과 무게를 고려해 1L로 결정" (the second numbered item),
#    insert new numbered paragraph "3. ".
# ---------------------------------------------------------------------
$pNumbered2 = $d.Paragraphs.Item(8)
$body2 = '<w:p><w:pPr><w:pStyle w:val="a4"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:leftChars="0"/></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">3. </w:t></w:r>' + `
    '</w:p>'
Insert-ParaAfter $pNumbered2 $body2

Write-Host "step2 done, count=$($d.Paragraphs.Count)"

# ---------------------------------------------------------------------
# 3) The "- LCD 모니터 사용법 확인" paragraph: drop the rPr from its pPr
#    (keep ind firstLine=195) and drop the bookmark (it moves to the new
#    last paragraph added in step 4).
# ---------------------------------------------------------------------
$pMonitor = $d.Paragraphs.Item(14)
$body3 = '<w:p><w:pPr><w:ind w:firstLine="195"/></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">LCD </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>모니터 사용법 확인</w:t></w:r>' + `
    '</w:p>'
Replace-ParaContent $pMonitor $body3

Write-Host "step3 done, count=$($d.Paragraphs.Count)"

# ---------------------------------------------------------------------
# 4) After the "모니터" paragraph, insert two new paragraphs:
#      "- 기타 재료 구매"
#      "- Garbage Collector 알고리즘 순서도 제작"  (carries the _GoBack bookmark)
# ---------------------------------------------------------------------
$pMonitor = $d.Paragraphs.Item(14)
$body4a = '<w:p><w:pPr><w:ind w:firstLine="195"/></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>- 기타 재료 구매</w:t></w:r>' + `
    '</w:p>'
Insert-ParaAfter $pMonitor $body4a

$pMaterials = $d.Paragraphs.Item(15)
$body4b = '<w:p><w:pPr><w:ind w:firstLine="195"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>- Garbage Collector 알고리즘 순서도 제작</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'
Insert-ParaAfter $pMaterials $body4b

Write-Host "step4 done, count=$($d.Paragraphs.Count)"

# ---------------------------------------------------------------------
# 5) The final (originally last) paragraph "- 기타 재료 구매" loses its run
#    (its text moved earlier in step 4); it keeps its original pPr but
#    becomes an otherwise-empty paragraph.
# ---------------------------------------------------------------------
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$body5 = '<w:p><w:pPr><w:ind w:firstLine="195"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>'
Replace-ParaContent $pLast $body5

Write-Host "step5 done, count=$($d.Paragraphs.Count)"

# ---------------------------------------------------------------------
# 6) Append the new "추가로 필요한 부품 및 재료" section after it.
# ---------------------------------------------------------------------
$pEmpty = $d.Paragraphs.Item($d.Paragraphs.Count)
$body6 = '<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>추가로 필요한 부품 및 재료:</w:t></w:r></w:p>'
Insert-ParaAfter $pEmpty $body6

$pHeader = $d.Paragraphs.Item($d.Paragraphs.Count)
$body7 = '<w:p><w:pPr><w:ind w:firstLine="195"/></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">- 바퀴 </w:t></w:r>' + `
    '<w:r><w:t>4</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>개</w:t></w:r>' + `
    '</w:p>'
Insert-ParaAfter $pHeader $body7

$pWheels = $d.Paragraphs.Item($d.Paragraphs.Count)
$body8 = '<w:p><w:pPr><w:ind w:firstLine="195"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">- 스텝 모터 </w:t></w:r>' + `
    '<w:r><w:t>4</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>개</w:t></w:r>' + `
    '</w:p>'
Insert-ParaAfter $pWheels $body8

$pStepMotor = $d.Paragraphs.Item($d.Paragraphs.Count)
$body9 = '<w:p><w:pPr><w:ind w:firstLine="195"/></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>서보</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 모터 </w:t></w:r>' + `
    '<w:r><w:t>1</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>개</w:t></w:r>' + `
    '</w:p>'
Insert-ParaAfter $pStepMotor $body9

$pServoMotor = $d.Paragraphs.Item($d.Paragraphs.Count)
$body10 = '<w:p><w:pPr><w:ind w:firstLine="195"/></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">- 초음파 센서 </w:t></w:r>' + `
    '<w:r><w:t>1</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>개</w:t></w:r>' + `
    '</w:p>'
Insert-ParaAfter $pServoMotor $body10

$pUltrasonic = $d.Paragraphs.Item($d.Paragraphs.Count)
$body11 = '<w:p><w:pPr><w:ind w:firstLine="195"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>- 리튬 이온 배터리 홀더</w:t></w:r>' + `
    '</w:p>'
Insert-ParaAfter $pUltrasonic $body11

Write-Host "step6 done, count=$($d.Paragraphs.Count)"
